$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $value) {
    $c = $ws.Range($cellref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.203.67"
Set-TextValue "E2" "  +0.31%  "
Set-TextValue "D3" "1.834.28"
Set-TextValue "E3" "  -0.18%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "242.76"
Set-TextValue "E5" "  +0.86%  "
Set-TextValue "D6" "0.6624"
Set-TextValue "E6" "  -2.68%  "
Set-TextValue "D7" "1.000"
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "0.07424"
Set-TextValue "E8" "  -0.54%  "
Set-TextValue "D9" "0.2939"
Set-TextValue "E9" "  -1.62%  "
Set-TextValue "D10" "22.89"
Set-TextValue "E10" "  -1.10%  "
Set-TextValue "D11" "0.07772"
Set-TextValue "E11" "  +1.42%  "
Set-TextValue "D12" "1.819.34"
Set-TextValue "E12" "  -0.77%  "
Set-TextValue "D13" "5.000"
Set-TextValue "E13" "  -0.52%  "
Set-TextValue "D14" "0.6676"
Set-TextValue "E14" "  -1.41%  "
Set-TextValue "D15" "82.94"
Set-TextValue "E15" "  -3.60%  "
Set-TextValue "D16" "6.118"
Set-TextValue "E16" "  -0.63%  "
Set-TextValue "D17" "0.000008390"
Set-TextValue "E17" "  +1.28%  "
Set-TextValue "D18" "29.189.75"
Set-TextValue "E18" "  +0.33%  "
Set-TextValue "D19" "2.075.53"
Set-TextValue "E19" "  +0.34%  "
Set-TextValue "D20" "228.51"
Set-TextValue "E20" "  +0.20%  "
Set-TextValue "D21" "12.48"
Set-TextValue "E21" "  -0.12%  "
Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  +0.14%  "
Set-TextValue "D23" "7.148"
Set-TextValue "E23" "  -2.73%  "
Set-TextValue "D24" "1.000"
Set-TextValue "D25" "159.27"
Set-TextValue "E25" "  -0.80%  "
Set-TextValue "B26" "Cosmos"
Set-TextValue "C26" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D26" "8.625"
Set-TextValue "E26" "  -0.95%  "
Set-TextValue "B27" "Stellar"
Set-TextValue "C27" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D27" "0.1401"
Set-TextValue "E27" "  -2.64%  "
Set-TextValue "D28" "18.01"
Set-TextValue "E28" "  -0.11%  "
Set-TextValue "D29" "1.518"
Set-TextValue "E29" "  +0.84%  "
Set-TextValue "D30" "4.119"
Set-TextValue "E30" "  -3.10%  "
Set-TextValue "D31" "4.051"
Set-TextValue "E31" "  -1.97%  "
Set-TextValue "E32" "  -0.32%  "
Set-TextValue "D33" "0.05281"
Set-TextValue "E33" "  -2.24%  "
Set-TextValue "D34" "1.866"
Set-TextValue "E34" "  +0.14%  "
Set-TextValue "D35" "0.7443"
Set-TextValue "E35" "  -0.75%  "
Set-TextValue "D36" "1.141"
Set-TextValue "E36" "  +1.02%  "
Set-TextValue "D37" "2.655"
Set-TextValue "E37" "  -0.99%  "
Set-TextValue "D38" "1.314.75"
Set-TextValue "E38" "  +1.02%  "
Set-TextValue "D39" "0.01796"
Set-TextValue "E39" "  -0.96%  "
Set-TextValue "D40" "2.740"
Set-TextValue "E40" "  +1.08%  "
Set-TextValue "D41" "0.9296"
Set-TextValue "E41" "  -0.81%  "
Set-TextValue "D42" "5.918"
Set-TextValue "E42" "  -2.72%  "
Set-TextValue "D43" "0.08433"
Set-TextValue "E43" "  +5.62%  "
Set-TextValue "D44" "0.9995"
Set-TextValue "E44" "  +0.11%  "
Set-TextValue "D45" "102.79"
Set-TextValue "E45" "  -1.82%  "
Set-TextValue "D46" "1.973.06"
Set-TextValue "E46" "  +0.19%  "
Set-TextValue "D47" "0.5144"
Set-TextValue "E47" "  -0.54%  "
Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.00000000120"
Set-TextValue "E48" "  -0.66%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.756"
Set-TextValue "E49" "  -0.45%  "
Set-TextValue "D50" "63.17"
Set-TextValue "E50" "  -1.12%  "
Set-TextValue "D51" "0.05876"
Set-TextValue "E51" "  -0.86%  "
